# Nexial desktop commands: add `clickScreen(button,modifiers,x,y)` and
# `mouseWheel(amount,modifiers,x,y)` to the alphabetically-sorted list of
# desktop commands kept on the hidden `#system` sheet (column G), which
# backs the `desktop` named range used for data-validation dropdowns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# The existing list lives in G2:G95. Two new entries are being spliced in:
#   - "clickScreen(button,modifiers,x,y)" belongs right after "clickRadio(name)" (row 41)
#   - "mouseWheel(amount,modifiers,x,y)" belongs right after "minimize()" (originally row 57)
# Shift everything below each insertion point down by one row (bottom-up,
# so we never overwrite a value before it has been copied).

# Step 1: make room for "mouseWheel" - push rows 57-95 down to 59-97
for ($r = 95; $r -ge 57; $r--) {
    $ws.Cells.Item($r + 2, 7).Value = $ws.Cells.Item($r, 7).Value2
}

# Step 2: make room for "clickScreen" - push rows 42-56 down to 43-57
for ($r = 56; $r -ge 42; $r--) {
    $ws.Cells.Item($r + 1, 7).Value = $ws.Cells.Item($r, 7).Value2
}

# Step 3: write the two new commands into the gaps just opened up
$ws.Cells.Item(42, 7).Value = "clickScreen(button,modifiers,x,y)"
$ws.Cells.Item(58, 7).Value = "mouseWheel(amount,modifiers,x,y)"

# Step 4: the `desktop` named range now spans two additional rows
$wb.Names.Item("desktop").RefersTo = "='#system'!`$G`$2:`$G`$97"
